$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "descr" column (C) values:
#  - French rows:  "Mot sur la liste noire" -> "Mot sur la liste de blocage"
#  - English rows: "Blacklisted Word" -> "Blocklisted Word"
# (French cells are updated first so the new shared strings are appended
# to the shared string table in the same order as the target workbook.)
$ws.Range("C6").Value = "Mot sur la liste de blocage"
$ws.Range("C7").Value = "Mot sur la liste de blocage"
$ws.Range("C2").Value = "Blocklisted Word"
$ws.Range("C3").Value = "Blocklisted Word"
$ws.Range("C4").Value = "Blocklisted Word"
$ws.Range("C5").Value = "Blocklisted Word"

# Update the active cell selection to C5 (was B5)
$ws.Range("C5").Select()
